$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "2023-09-02 11:18:06"
$ws.Range("B2").Value = "Entrada"
$ws.Range("C2").Value = "Vamos ver"
$ws.Range("D2").Value = 11
$ws.Range("E2").Value = "Novo"

# Update row 3
$ws.Range("A3").Value = "2023-09-02 11:49:24"
$ws.Range("B3").Value = "Entrada"
$ws.Range("C3").Value = "Vamos ver"
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = "vamnos"

# Remove rows 4-7 (no longer present in the data)
$ws.Range("A4:E7").Delete()
